$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# --- Sheet "VENTAS POR GRUPO": zero out cells affected by corrections ---
$ws1.Range("H20").Value = 0
$ws1.Range("I20").Value = 0
$ws1.Range("M20").Value = 0
$ws1.Range("L27").Value = 0
$ws1.Range("L33").Value = 0
$ws1.Range("L42").Value = 0
$ws1.Range("L43").Value = 0
$ws1.Range("E45").Value = 0
$ws1.Range("L45").Value = 0
$ws1.Range("E54").Value = 0
$ws1.Range("M54").Value = 0

# --- Sheet "VENTAS POR GRUPO": update "x de 56" progress counters in row 58 ---
$ws1.Range("E58").Value = "0 de 56"
$ws1.Range("H58").Value = "0 de 56"
$ws1.Range("I58").Value = "0 de 56"
$ws1.Range("L58").Value = "0 de 56"
$ws1.Range("M58").Value = "0 de 56"

# --- Sheet "VENTA MENSUAL": roll month headers forward one column ---
$ws2.Range("C1").Value = "septiembre"
$ws2.Range("D1").Value = "octubre"
$ws2.Range("E1").Value = "noviembre"
$ws2.Range("F1").Value = "diciembre"

# --- Sheet "VENTA MENSUAL": column widths follow the shifted month columns ---
$ws2.Columns.Item(3).ColumnWidth = 15.17
$ws2.Columns.Item(4).ColumnWidth = 12.17
$ws2.Columns.Item(5).ColumnWidth = 14.17

# --- Sheet "VENTA MENSUAL": shift monthly sales figures one column to the right ---
$ws2.Range("C5").Value = 0
$ws2.Range("C6").Value = 0
$ws2.Range("C9").Value = 3864.45
$ws2.Range("D9").Value = -22.29
$ws2.Range("E9").Value = 0
$ws2.Range("C10").Value = 236.29
$ws2.Range("D10").Value = 0
$ws2.Range("C12").Value = 448.77
$ws2.Range("D12").Value = 0
$ws2.Range("C13").Value = 236.29
$ws2.Range("D13").Value = 0
$ws2.Range("C14").Value = 0
$ws2.Range("C16").Value = 508.48
$ws2.Range("D16").Value = 0
$ws2.Range("E20").Value = 3669.5
$ws2.Range("F20").Value = 0
$ws2.Range("C21").Value = 0
$ws2.Range("D23").Value = 56.02
$ws2.Range("E23").Value = 0
$ws2.Range("E27").Value = 393.88
$ws2.Range("F27").Value = 0
$ws2.Range("C30").Value = 430.11
$ws2.Range("D30").Value = 0
$ws2.Range("E33").Value = 216.76
$ws2.Range("F33").Value = 0
$ws2.Range("C34").Value = 551.71
$ws2.Range("D34").Value = 0
$ws2.Range("D35").Value = -166.48
$ws2.Range("E35").Value = 0
$ws2.Range("C39").Value = 0
$ws2.Range("C40").Value = 0
$ws2.Range("E42").Value = 450.22
$ws2.Range("F42").Value = 0
$ws2.Range("C43").Value = 115.52
$ws2.Range("D43").Value = 179.12
$ws2.Range("E43").Value = 517.06
$ws2.Range("F43").Value = 0
$ws2.Range("E45").Value = 503.15
$ws2.Range("F45").Value = 0
$ws2.Range("C47").Value = 44.79
$ws2.Range("D47").Value = 0
$ws2.Range("D48").Value = 194.16
$ws2.Range("E48").Value = 0
$ws2.Range("C51").Value = 5858.53
$ws2.Range("D51").Value = 0
$ws2.Range("D52").Value = 334.37
$ws2.Range("E52").Value = 0
$ws2.Range("C53").Value = -10.44
$ws2.Range("D53").Value = 0
$ws2.Range("C54").Value = 522.82
$ws2.Range("D54").Value = 0
$ws2.Range("E54").Value = 944.42
$ws2.Range("F54").Value = 0
$ws2.Range("C55").Value = 165.83
$ws2.Range("D55").Value = 0
$ws2.Range("C57").Value = 438.86
$ws2.Range("D57").Value = 0
$ws2.Range("C58").Value = 13412.01
$ws2.Range("D58").Value = 574.9
$ws2.Range("E58").Value = 6694.99
$ws2.Range("F58").Value = 0

Write-Host "Edit applied."
